# Split the run containing
#   "am STERBEDATUM im Alter von ALTER Jahren, hier in STERBEORT."
# into three runs (same Arial/NewCenturySchlbk-Roman/26pt formatting),
# inserting a new "LEBENS" run right before "ALTER", so the merge-field
# placeholder "ALTER" becomes "LEBENSALTER":
#   "am STERBEDATUM im Alter von " | "LEBENS" | "ALTER Jahren, hier in STERBEORT."

$d = $word.ActiveDocument

$oldSentence = "am STERBEDATUM im Alter von ALTER Jahren, hier in STERBEORT."

$search = $d.Content
$found = $search.Find.Execute($oldSentence, $true, $false, $false, $false, `
                               $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target sentence '$oldSentence'"
}

# Re-materialize a plain Range over the found span (InsertXML needs a
# freshly bound Range, not the Find-owning one, to replace-in-place
# instead of inserting at its start).
$target = $d.Range($search.Start, $search.End)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Arial" w:eastAsia="NewCenturySchlbk-Roman" w:hAnsi="Arial" w:cs="Arial"/>
<w:sz w:val="26"/>
<w:szCs w:val="26"/>
</w:rPr>
<w:t xml:space="preserve">am STERBEDATUM im Alter von </w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Arial" w:eastAsia="NewCenturySchlbk-Roman" w:hAnsi="Arial" w:cs="Arial"/>
<w:sz w:val="26"/>
<w:szCs w:val="26"/>
</w:rPr>
<w:t>LEBENS</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Arial" w:eastAsia="NewCenturySchlbk-Roman" w:hAnsi="Arial" w:cs="Arial"/>
<w:sz w:val="26"/>
<w:szCs w:val="26"/>
</w:rPr>
<w:t>ALTER Jahren, hier in STERBEORT.</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target.InsertXML($xml)

Write-Output "Split 'am STERBEDATUM...' run into 3 runs with LEBENS inserted before ALTER."
